$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C50").Value = "Terminé"
$ws.Range("E50").Value = 4
$ws.Range("C51").Value = "En cours"

$ws.Range("J53").Select()
